$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E ("Komentarz" for the third pair Łukasz Napora) gets the same
# comment "Refaktoryzacja kodu" added for weeks 6-13 (rows 10-17).
$ws.Range("E10").Value = "Refaktoryzacja kodu"
$ws.Range("E11").Value = "Refaktoryzacja kodu"
$ws.Range("E12").Value = "Refaktoryzacja kodu"
$ws.Range("E13").Value = "Refaktoryzacja kodu"
$ws.Range("E14").Value = "Refaktoryzacja kodu"
$ws.Range("E15").Value = "Refaktoryzacja kodu"
$ws.Range("E16").Value = "Refaktoryzacja kodu"
$ws.Range("E17").Value = "Refaktoryzacja kodu"

# Column C ("Komentarz" for the first pair Michał Mierzyński) gets new
# per-row comments for weeks 7-13 (rows 11-17).
$ws.Range("C11").Value = "Implementacja metod Merge i Divide. Integracja z komponentami"
$ws.Range("C12").Value = "Implementacja metod Merge i Divide. Integracja z komponentami"
$ws.Range("C13").Value = "Implementacja metod Merge i Divide. Integracja z komponentami"
$ws.Range("C14").Value = "Popawki wydajnościowe w komunikacji oraz w metodach wtyczki"
$ws.Range("C15").Value = "Implementacja metod Merge i Divide. Integracja z komponentami innych grup"
$ws.Range("C16").Value = "Implementacja metod Merge i Divide. Integracja z komponentami innych grup"
$ws.Range("C17").Value = "Implementacja metod Merge i Divide. Integracja z komponentami innych grup"

# Move the active selection to C17 (was G10), and the view no longer pins
# topLeftCell to A8.
$ws.Range("C17").Select()
